{"js": "// Update the date and the division problems for this worksheet.\nconst replacements = [\n    [\"2024-08-10 Saturday\", \"2024-08-11 Sunday\"],\n    [\"752\u00f73=\", \"287\u00f78=\"],\n    [\"346\u00f74=\", \"448\u00f79=\"],\n    [\"722\u00f78=\", \"674\u00f75=\"],\n    [\"587\u00f78=\", \"432\u00f76=\"],\n    [\"832\u00f79=\", \"540\u00f74=\"],\n    [\"159\u00f72=\", \"249\u00f78=\"],\n    [\"249\u00f79=\", \"882\u00f72=\"],\n    [\"865\u00f79=\", \"210\u00f72=\"],\n    [\"629\u00f78=\", \"745\u00f76=\"],\n    [\"880\u00f74=\", \"305\u00f72=\"],\n    [\"418\u00f73=\", \"123\u00f79=\"],\n    [\"978\u00f75=\", \"776\u00f73=\"],\n    [\"375\u00f77=\", \"903\u00f79=\"],\n    [\"633\u00f77=\", \"289\u00f78=\"],\n    [\"706\u00f74=\", \"670\u00f79=\"],\n    [\"710\u00f78=\", \"411\u00f77=\"],\n    [\"494\u00f75=\", \"540\u00f78=\"],\n    [\"509\u00f75=\", \"413\u00f74=\"],\n    [\"385\u00f75=\", \"836\u00f72=\"],\n    [\"520\u00f77=\", \"608\u00f76=\"],\n    [\"628\u00f74=\", \"155\u00f72=\"],\n    [\"494\u00f76=\", \"365\u00f74=\"],\n    [\"844\u00f79=\", \"815\u00f74=\"],\n    [\"429\u00f75=\", \"606\u00f75=\"],\n    [\"250\u00f79=\", \"949\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n    const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (const range of results.items) {\n        range.insertText(newText, \"Replace\");\n    }\n    await context.sync();\n}\n", "ps1": "# Update the date and the division problems for this worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = '2024-08-10 Saturday'; New = '2024-08-11 Sunday'},\n    @{Old = '752\u00f73='; New = '287\u00f78='},\n    @{Old = '346\u00f74='; New = '448\u00f79='},\n    @{Old = '722\u00f78='; New = '674\u00f75='},\n    @{Old = '587\u00f78='; New = '432\u00f76='},\n    @{Old = '832\u00f79='; New = '540\u00f74='},\n    @{Old = '159\u00f72='; New = '249\u00f78='},\n    @{Old = '249\u00f79='; New = '882\u00f72='},\n    @{Old = '865\u00f79='; New = '210\u00f72='},\n    @{Old = '629\u00f78='; New = '745\u00f76='},\n    @{Old = '880\u00f74='; New = '305\u00f72='},\n    @{Old = '418\u00f73='; New = '123\u00f79='},\n    @{Old = '978\u00f75='; New = '776\u00f73='},\n    @{Old = '375\u00f77='; New = '903\u00f79='},\n    @{Old = '633\u00f77='; New = '289\u00f78='},\n    @{Old = '706\u00f74='; New = '670\u00f79='},\n    @{Old = '710\u00f78='; New = '411\u00f77='},\n    @{Old = '494\u00f75='; New = '540\u00f78='},\n    @{Old = '509\u00f75='; New = '413\u00f74='},\n    @{Old = '385\u00f75='; New = '836\u00f72='},\n    @{Old = '520\u00f77='; New = '608\u00f76='},\n    @{Old = '628\u00f74='; New = '155\u00f72='},\n    @{Old = '494\u00f76='; New = '365\u00f74='},\n    @{Old = '844\u00f79='; New = '815\u00f74='},\n    @{Old = '429\u00f75='; New = '606\u00f75='},\n    @{Old = '250\u00f79='; New = '949\u00f74='}\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $pair.Old\n    $rng.Find.Replacement.Text = $pair.New\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.Wrap = 1  # wdFindContinue\n    $rng.Find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)  # wdReplaceAll = 2\n}\n"}
